$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for "2022-Q3" right above the
#    existing "2022-Q2" row, pushing everything else down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A2").EntireRow.Insert()
$total.Range("A2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 1.96

# Pick up the "A column" style (bold-ish index used by the rest of the
# column) from the row just below, which already has it, instead of
# hand-building a style.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)   # xlPasteFormats
$total.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2) New "2022-Q3" worksheet, positioned right after "总计" (i.e. right
#    before "2022-Q2"). Duplicating the existing "2022-Q2" sheet gives us
#    the same sheetPr/pageMargins/header-style scaffolding, then we trim
#    the extra row and overwrite the data values.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The template sheet has 8 data rows (2..9); the new sheet only needs 7 (2..8)
$q3.Range("A9:H9").Delete(-4162)   # xlShiftUp

$rows = @(
  @(0, "010967", "博道嘉丰混合A",                     "13.90", "89.51", "4.21", "0.5852", 6),
  @(1, "010147", "博道嘉兴一年持有期混合",               "11.13", "90.11", "4.69", "0.5220", 7),
  @(2, "008467", "博道嘉瑞混合A",                     "13.34", "77.63", "2.63", "0.3508", 9),
  @(3, "008208", "博道嘉泰回报混合",                   "12.96", "77.67", "2.67", "0.3460", 9),
  @(4, "010968", "博道嘉丰混合C",                     "3.17",  "89.51", "4.21", "0.1335", 6),
  @(5, "008468", "博道嘉瑞混合C",                     "0.72",  "77.63", "2.63", "0.0189", 9),
  @(6, "562530", "华夏中证智选1000价值稳健策略ETF",       "0.54",  "94.32", "0.95", "0.0051", 4)
)

$textCols = @("B", "C", "D", "E", "F", "G")

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $i + 2
  $row = $rows[$i]

  $q3.Range("A$r").Value = $row[0]

  # Force text storage for B..G so numeric-looking strings ("010967",
  # "13.90", ...) are kept verbatim instead of being auto-coerced to
  # numbers by the usual Excel type-inference on Value assignment.
  foreach ($col in $textCols) {
    $q3.Range("$col$r").NumberFormat = "@"
  }

  $q3.Range("B$r").Value = $row[1]
  $q3.Range("C$r").Value = $row[2]
  $q3.Range("D$r").Value = $row[3]
  $q3.Range("E$r").Value = $row[4]
  $q3.Range("F$r").Value = $row[5]
  $q3.Range("G$r").Value = $row[6]

  # Drop the now-unneeded "@" number format so the cells end up with no
  # special style, matching the rest of the sheet.
  foreach ($col in $textCols) {
    $q3.Range("$col$r").ClearFormats()
  }

  $q3.Range("H$r").Value = $row[7]
}

Write-Output "2022-Q3 sheet inserted and 总计 updated"
